# Apply cryptos list price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.803.76"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.925.61"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'375.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'100.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "'35.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "3.389.85"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'7.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "'11.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +54.51%  "
$ws.Range("D17").Value = "2.921.68"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "50.790.74"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "'3.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.98%  "
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'68.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'264.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "'3.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.11%  "
$ws.Range("D26").Value = "'8.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'7.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'25.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'9.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "'50.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "'32.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'3.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'16.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'2.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").Value = "'118.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").Value = "'21.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").Value = "'3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.41%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'0.265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").Value = "1.988.49"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "'0.0322"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("E51").Value = "  +1.99%  "
